$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.967.62"
$ws.Range("E2").Value = "  +2.19%  "
$ws.Range("D3").Value = "3.172.64"
$ws.Range("E3").Value = "  +4.37%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.71"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +4.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.99"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +7.37%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.174.32"
$ws.Range("E8").Value = "  +4.47%  "
$ws.Range("E9").Value = "  +3.16%  "
$ws.Range("E10").Value = "  +6.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.24"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.503"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +3.26%  "
$ws.Range("E13").Value = "  +18.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.71"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +6.23%  "
$ws.Range("D15").Value = "3.694.48"
$ws.Range("E15").Value = "  +4.49%  "
$ws.Range("D16").Value = "65.036.41"
$ws.Range("E16").Value = "  +2.20%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.18"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +6.13%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.169.44"
$ws.Range("E18").Value = "  +4.16%  "
$ws.Range("E19").Value = "  +1.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "514.22"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +8.40%  "
$ws.Range("E21").Value = "  +6.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.729"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +7.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.30"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +4.74%  "
$ws.Range("E24").Value = "  +4.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.40"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.45%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.02"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +11.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.94"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +5.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.19"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +7.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "27.88"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +6.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.76"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +13.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.11%  "
$ws.Range("E33").Value = "  +5.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.36"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +10.89%  "
$ws.Range("E35").Value = "  +6.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.77"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0895"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +10.58%  "
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.16"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +13.71%  "
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "475.07"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +7.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0420"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +3.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.66"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +4.98%  "
$ws.Range("D42").Value = "3.065.31"
$ws.Range("E42").Value = "  +1.76%  "
$ws.Range("E43").Value = "  +1.73%  "
$ws.Range("E44").Value = "  +6.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.42"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +8.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.09"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.96%  "
$ws.Range("D47").Value = "0.0₃0612"
$ws.Range("E47").Value = "  +19.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.115"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.15%  "
$ws.Range("E50").Value = "  +8.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.51"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.31%  "
